# Automatische test-sync: 2025-08-13 22:16:50
# Appends the new mail-log entry to the "Logs" sheet and refreshes the
# "Dashboard" summary count to match.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

# Find the next empty row right after the current used range.
$newRow = $logs.Cells.Item($logs.UsedRange.Rows.Count + 1, 1).Row

$logs.Cells.Item($newRow, 1).Value = "Demo inplannen"
$logs.Cells.Item($newRow, 2).Value = "klantenservice@testbedrijf123.nl"
$logs.Cells.Item($newRow, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-13 22:16:16"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 23

# Extend the conditional-formatting ranges so they keep covering the full
# data set now that a new row has been appended.
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$($col)2:$($col)23")
    $newRange = $logs.Range("$($col)2:$($col)24")
    $fcs = $oldRange.FormatConditions
    if ($fcs.Count -gt 0) {
        $fcs.Item(1).ModifyAppliesToRange($newRange)
    }
}
